$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Day 4")
$ws.Activate()
try {
  $excel.Goto($ws.Range("A225"), $true)
  Write-Output "Goto ok"
} catch {
  Write-Output "Goto failed: $_"
}
